# "Reload OK...need some improvement"
#
# Insert a new task row ("Reload" / 重新啟動) into the Task sheet just above
# the existing "ShuangShiyiActivityReward" row (row 23), pushing the rest of
# the table down by one row. Mark the new task's Priority cell (column D)
# as "HIGH" with a red highlight fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task")

# Push rows 23..32 down to 24..33, duplicating row 22's formatting into the
# freshly inserted row 23 (matches native Excel Insert-row behaviour).
$ws.Rows.Item(23).Insert()

# New task data.
$ws.Cells.Item(23, 1).Value = "Reload"
$ws.Cells.Item(23, 2).Value = "重新啟動"
$ws.Cells.Item(23, 3).Value = "Yes"
$ws.Cells.Item(23, 4).Value = "HIGH"

# Highlight the new Priority cell in red.
$ws.Cells.Item(23, 4).Interior.Color = 255

# Restore the "case TaskId.X:" helper formula for the new row (column O).
$ws.Cells.Item(23, 15).Formula = "=""case TaskId.""&A23&"":"""

# Leave the selection where the author left off editing.
$ws.Range("D23").Select() | Out-Null
